# Insert a new data row at row 298 (shifts existing rows 298..402 down to 299..403)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(298).Insert()

# Populate the newly inserted row 298 with the new record's data.
$ws.Range("A298").Value = 3
$ws.Range("B298").Value = "Femacal de La Calera"
$ws.Range("C298").Value = "Coquimbo"
$ws.Range("D298").Value = 44627
$ws.Range("D298").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E298").Value = 5
$ws.Range("F298").Value = 100112021
$ws.Range("G298").Value = "Ají"
$ws.Range("H298").Value = "Americana (o)"
$ws.Range("I298").Value = "Primera"
$ws.Range("J298").Value = 85
$ws.Range("K298").Value = 15000
$ws.Range("L298").Value = 16000
$ws.Range("M298").Value = 15529
$ws.Range("N298").Value = '$/caja 15 kilos'
$ws.Range("O298").Value = "Provincia de Quillota"
$ws.Range("P298").Value = 1035
$ws.Range("Q298").Value = 15
$ws.Range("R298").Value = "Hortaliza"
